$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate column C with country names (mirrors a VLOOKUP against a reference table, pasted as values)
$ws.Range("C1").Value2 = "Pays_lbl"
$ws.Range("C2").Value2 = "Kazakhstan"
$ws.Range("C3").Value2 = "Espagne"
$ws.Range("C4").Value2 = "France"
$ws.Range("C5").Value2 = "France"
$ws.Range("C6").Value2 = "Pays-Bas"
$ws.Range("C7").Value2 = "Pays-Bas"
$ws.Range("C8").Value2 = "France"
$ws.Range("C9").Value2 = "France"
$ws.Range("C10").Value2 = "Suisse"
$ws.Range("C11").Value2 = "Belgique"
$ws.Range("C12").Value2 = "Italie"
$ws.Range("C13").Value2 = "Belgique"
$ws.Range("C14").Value2 = "Estonie"
$ws.Range("C15").Value2 = "France"
$ws.Range("C16").Value2 = "Ukraine"
$ws.Range("C17").Value2 = "Espagne"
$ws.Range("C18").Value2 = "Pologne"
$ws.Range("C19").Value2 = "Norvège"
$ws.Range("C20").Value2 = "Espagne"
$ws.Range("C21").Value2 = "Lituanie"
$ws.Range("C22").Value2 = "Royaume-Uni"
$ws.Range("C23").Value2 = "#N/A"
$ws.Range("C24").Value2 = "Portugal"
$ws.Range("C25").Value2 = "France"
$ws.Range("C26").Value2 = "Espagne"
$ws.Range("C27").Value2 = "Australie"
$ws.Range("C28").Value2 = "Italie"
$ws.Range("C29").Value2 = "#N/A"
$ws.Range("C30").Value2 = "France"
$ws.Range("C31").Value2 = "Danemark"
$ws.Range("C32").Value2 = "Colombie"
$ws.Range("C33").Value2 = "Italie"
$ws.Range("C34").Value2 = "Allemagne"
$ws.Range("C35").Value2 = "Danemark"
$ws.Range("C36").Value2 = "Belgique"
$ws.Range("C37").Value2 = "France"
$ws.Range("C38").Value2 = "Espagne"
$ws.Range("C39").Value2 = "Italie"
$ws.Range("C40").Value2 = "Italie"
$ws.Range("C41").Value2 = "Espagne"
$ws.Range("C42").Value2 = "Autriche"
$ws.Range("C43").Value2 = "Italie"
$ws.Range("C44").Value2 = "Pays-Bas"
$ws.Range("C45").Value2 = "Pays-Bas"
$ws.Range("C46").Value2 = "Russie"
$ws.Range("C47").Value2 = "Biélorussie"
$ws.Range("C48").Value2 = "Allemagne"
$ws.Range("C49").Value2 = "Espagne"
$ws.Range("C50").Value2 = "France"
$ws.Range("C51").Value2 = "France"
$ws.Range("C52").Value2 = "Japon"
$ws.Range("C53").Value2 = "Espagne"
$ws.Range("C54").Value2 = "Italie"
$ws.Range("C55").Value2 = "Nouvelle Zélande"
$ws.Range("C56").Value2 = "Belgique"
$ws.Range("C57").Value2 = "Italie"
$ws.Range("C58").Value2 = "Slovaquie"
$ws.Range("C59").Value2 = "Italie"
$ws.Range("C60").Value2 = "Allemagne"
$ws.Range("C61").Value2 = "Italie"
$ws.Range("C62").Value2 = "Italie"
$ws.Range("C63").Value2 = "Espagne"
$ws.Range("C64").Value2 = "Espagne"
$ws.Range("C65").Value2 = "France"
$ws.Range("C66").Value2 = "États-Unis (USA)"
$ws.Range("C67").Value2 = "Danemark"
$ws.Range("C68").Value2 = "Slovénie"
$ws.Range("C69").Value2 = "Czechia"
$ws.Range("C70").Value2 = "Allemagne"
$ws.Range("C71").Value2 = "France"
$ws.Range("C72").Value2 = "France"
$ws.Range("C73").Value2 = "Argentine"
$ws.Range("C74").Value2 = "#N/A"
$ws.Range("C75").Value2 = "Canada"
$ws.Range("C76").Value2 = "Mexique"
$ws.Range("C77").Value2 = "Espagne"
$ws.Range("C78").Value2 = "Irlande"
$ws.Range("C79").Value2 = "France"
$ws.Range("C80").Value2 = "France"
$ws.Range("C81").Value2 = "Pays-Bas"
$ws.Range("C82").Value2 = "Pays-Bas"
$ws.Range("C83").Value2 = "Éthiopie"
$ws.Range("C84").Value2 = "Allemagne"
$ws.Range("C85").Value2 = "Pays-Bas"
$ws.Range("C86").Value2 = "#N/A"
$ws.Range("C87").Value2 = "Brésil"
$ws.Range("C88").Value2 = "Suisse"
$ws.Range("C89").Value2 = "Portugal"
$ws.Range("C90").Value2 = "Allemagne"
$ws.Range("C91").Value2 = "Luxembourg"
$ws.Range("C92").Value2 = "Érythrée"
$ws.Range("C93").Value2 = "Algérie"
$ws.Range("C94").Value2 = "Finlande"
$ws.Range("C95").Value2 = "Allemagne"
$ws.Range("C96").Value2 = "Taiwan"
$ws.Range("C97").Value2 = "Allemagne"
$ws.Range("C98").Value2 = "Danemark"
$ws.Range("C99").Value2 = "France"
$ws.Range("C100").Value2 = "Italie"
$ws.Range("C101").Value2 = "#N/A"
$ws.Range("C102").Value2 = "Turquie"
$ws.Range("C103").Value2 = "Espagne"
$ws.Range("C104").Value2 = "Belgique"
$ws.Range("C105").Value2 = "Danemark"
$ws.Range("C106").Value2 = "#N/A"
$ws.Range("C107").Value2 = "Allemagne"
$ws.Range("C108").Value2 = "Espagne"
$ws.Range("C109").Value2 = "Azerbaïdjan"
$ws.Range("C110").Value2 = "Italie"
$ws.Range("C111").Value2 = "#N/A"
$ws.Range("C112").Value2 = "Pays-Bas"
$ws.Range("C113").Value2 = "Pays-Bas"
$ws.Range("C114").Value2 = "France"
$ws.Range("C115").Value2 = "Italie"
$ws.Range("C116").Value2 = "Italie"
$ws.Range("C117").Value2 = "France"
$ws.Range("C118").Value2 = "France"
$ws.Range("C119").Value2 = "Vénézuela"
$ws.Range("C120").Value2 = "Espagne"
$ws.Range("C121").Value2 = "France"
$ws.Range("C122").Value2 = "Italie"
$ws.Range("C123").Value2 = "Allemagne"
$ws.Range("C124").Value2 = "#N/A"
$ws.Range("C125").Value2 = "Équateur"
$ws.Range("C126").Value2 = "Italie"
$ws.Range("C127").Value2 = "Israël"
$ws.Range("C128").Value2 = "Émirats Arabes Unis"
$ws.Range("C129").Value2 = "Maroc"
$ws.Range("C130").Value2 = "Allemagne"
$ws.Range("C131").Value2 = "Bulgarie"
$ws.Range("C132").Value2 = "#N/A"
$ws.Range("C133").Value2 = "Kuwait"
$ws.Range("C134").Value2 = "Rwanda"
$ws.Range("C135").Value2 = "Philippines"
$ws.Range("C136").Value2 = "Espagne"
$ws.Range("C137").Value2 = "#N/A"
$ws.Range("C138").Value2 = "Monaco"
$ws.Range("C139").Value2 = "Hongrie"
$ws.Range("C140").Value2 = "#N/A"
$ws.Range("C141").Value2 = "Suisse"
$ws.Range("C142").Value2 = "Panama"
$ws.Range("C143").Value2 = "Bahamas"
$ws.Range("C144").Value2 = "Ouzbékistan"
$ws.Range("C145").Value2 = "Allemagne"
$ws.Range("C146").Value2 = "France"
$ws.Range("C147").Value2 = "Italie"
$ws.Range("C148").Value2 = "Albanie"
$ws.Range("C149").Value2 = "Espagne"
$ws.Range("C150").Value2 = "Bahreïn"
$ws.Range("C151").Value2 = "Italie"
$ws.Range("C152").Value2 = "Puerto Rico"
$ws.Range("C153").Value2 = "Allemagne"
$ws.Range("C154").Value2 = "Italie"

# Column C width (bestFit-like)
$ws.Columns.Item(3).ColumnWidth = 16.35

# Re-apply AutoFilter over the extended range
$ws.Range("A1:C253").AutoFilter(1)

# Update the _FilterDatabase defined name to match the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Feuil1!_FilterDatabase") {
        $n.RefersTo = "=Feuil1!`$A`$1:`$C`$253"
    }
}

# Sheet2 page setup
$ws2 = $wb.Worksheets.Item(2)
$ps = $ws2.PageSetup
$ps.PaperSize = 9
$ps.Orientation = 1

Write-Output "done"
